$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.067.28"
$ws.Range("E2").Value = "  -2.05%  "

$ws.Range("D3").Value = "1.799.02"
$ws.Range("E3").Value = "  -2.62%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'307.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.63%  "

$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'0.4208"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.54%  "

$ws.Range("E8").Value = "  -2.81%  "

$ws.Range("D9").Value = "'0.07110"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.17%  "

$ws.Range("D10").Value = "'0.8440"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.69%  "

$ws.Range("D11").Value = "'20.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.03%  "

$ws.Range("D12").Value = "1.808.55"
$ws.Range("E12").Value = "  -4.75%  "

$ws.Range("D13").Value = "'5.292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").Value = "'6.367"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.62%  "

$ws.Range("D15").Value = "'0.06767"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "'80.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("D18").Value = "'0.000008724"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.91%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").Value = "'15.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.81%  "

$ws.Range("D21").Value = "27.064.60"
$ws.Range("E21").Value = "  -2.68%  "

$ws.Range("D22").Value = "'5.062"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").Value = "'11.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "2.017.34"
$ws.Range("E24").Value = "  -4.67%  "

$ws.Range("D25").Value = "'1.925"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.20%  "

$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("D27").Value = "'18.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.95%  "

$ws.Range("D28").Value = "'5.014"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.05%  "

$ws.Range("D29").Value = "'113.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.37%  "

$ws.Range("E30").Value = "  -12.55%  "

$ws.Range("D31").Value = "'0.09011"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").Value = "'0.7236"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.00%  "

$ws.Range("D33").Value = "'2.871"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.77%  "

$ws.Range("D34").Value = "'4.328"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").Value = "'1.090"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.03%  "

$ws.Range("D36").Value = "'1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").Value = "'1.081"
$ws.Range("D37").Style = "Normal"

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01906"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.65%  "

$ws.Range("D40").Value = "'0.1628"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.90%  "

$ws.Range("D41").Value = "'0.4963"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.36%  "

$ws.Range("D42").Value = "'2.612"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.60%  "

$ws.Range("D43").Value = "'8.035"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.25%  "

$ws.Range("D44").Value = "'5.896"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -13.08%  "

$ws.Range("D45").Value = "'105.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "'0.06296"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.88%  "

$ws.Range("D49").Value = "'0.4524"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.73%  "

$ws.Range("E50").Value = "  -4.29%  "

$ws.Range("D51").Value = "'1.711"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.79%  "
